$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Debit")
$ws.Range("B2").Value = "Fri Aug 22 18:01:57 EDT 2025"
$ws.Range("B3").Value = "Fri Aug 22 18:02:02 EDT 2025"
$ws.Range("B4").Value = "Fri Aug 22 18:02:06 EDT 2025"
$ws.Range("B5").Value = "Fri Aug 22 18:02:11 EDT 2025"
$ws.Range("B6").Value = "Fri Aug 22 18:02:15 EDT 2025"
$ws.Range("B7").Value = "Fri Aug 22 18:02:20 EDT 2025"
$ws.Range("B8").Value = "Fri Aug 22 18:02:24 EDT 2025"

$ws = $wb.Worksheets.Item("Debit-ZeroDollar")
$ws.Range("B2").Value = "Fri Aug 22 18:02:29 EDT 2025"
$ws.Range("B3").Value = "Fri Aug 22 18:02:33 EDT 2025"
$ws.Range("B4").Value = "Fri Aug 22 18:02:38 EDT 2025"
$ws.Range("B5").Value = "Fri Aug 22 18:02:42 EDT 2025"
$ws.Range("B6").Value = "Fri Aug 22 18:02:47 EDT 2025"
$ws.Range("B7").Value = "Fri Aug 22 18:02:51 EDT 2025"
$ws.Range("B8").Value = "Fri Aug 22 18:02:56 EDT 2025"

$ws = $wb.Worksheets.Item("Debit-Void")
$ws.Range("B2").Value = "Fri Aug 22 18:03:01 EDT 2025"
$ws.Range("B3").Value = "Fri Aug 22 18:03:09 EDT 2025"
$ws.Range("B4").Value = "Fri Aug 22 18:03:17 EDT 2025"
$ws.Range("B5").Value = "Fri Aug 22 18:03:25 EDT 2025"
$ws.Range("B6").Value = "Fri Aug 22 18:03:34 EDT 2025"
$ws.Range("B7").Value = "Fri Aug 22 18:03:42 EDT 2025"
$ws.Range("B8").Value = "Fri Aug 22 18:03:50 EDT 2025"

$ws = $wb.Worksheets.Item("Debit-Credit")
$ws.Range("B2").Value = "Fri Aug 22 18:03:58 EDT 2025"
$ws.Range("B3").Value = "Fri Aug 22 18:04:07 EDT 2025"
$ws.Range("B4").Value = "Fri Aug 22 18:04:16 EDT 2025"
$ws.Range("B5").Value = "Fri Aug 22 18:04:25 EDT 2025"
$ws.Range("B6").Value = "Fri Aug 22 18:04:33 EDT 2025"
$ws.Range("B7").Value = "Fri Aug 22 18:04:41 EDT 2025"
$ws.Range("B8").Value = "Fri Aug 22 18:04:50 EDT 2025"

$ws = $wb.Worksheets.Item("Debit-Credit-Void")
$ws.Range("B2").Value = "Fri Aug 22 18:04:59 EDT 2025"
$ws.Range("B3").Value = "Fri Aug 22 18:05:12 EDT 2025"
$ws.Range("B4").Value = "Fri Aug 22 18:05:29 EDT 2025"
$ws.Range("B5").Value = "Fri Aug 22 18:05:41 EDT 2025"
$ws.Range("B6").Value = "Fri Aug 22 18:05:54 EDT 2025"
$ws.Range("B7").Value = "Fri Aug 22 18:06:06 EDT 2025"
$ws.Range("B8").Value = "Fri Aug 22 18:06:20 EDT 2025"

$ws = $wb.Worksheets.Item("Debit-MRF")
$ws.Range("B2").Value = "Thu Aug 21 02:25:52 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 21 02:25:55 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 21 02:25:57 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 21 02:26:00 EDT 2025"
$ws.Range("B6").Value = "Thu Aug 21 02:26:03 EDT 2025"
$ws.Range("B7").Value = "Thu Aug 21 02:26:06 EDT 2025"
$ws.Range("B8").Value = "Thu Aug 21 02:26:09 EDT 2025"
$ws.Range("B9").Value = "Thu Aug 21 02:26:11 EDT 2025"
$ws.Range("B10").Value = "Thu Aug 21 02:26:15 EDT 2025"
$ws.Range("B11").Value = "Thu Aug 21 02:26:18 EDT 2025"
$ws.Range("B12").Value = "Thu Aug 21 02:26:21 EDT 2025"
$ws.Range("B13").Value = "Thu Aug 21 02:26:24 EDT 2025"
$ws.Range("B14").Value = "Thu Aug 21 02:26:27 EDT 2025"
$ws.Range("B15").Value = "Thu Aug 21 02:26:30 EDT 2025"
$ws.Range("B16").Value = "Thu Aug 21 02:26:33 EDT 2025"
$ws.Range("B17").Value = "Thu Aug 21 02:26:37 EDT 2025"
$ws.Range("B18").Value = "Thu Aug 21 02:26:40 EDT 2025"
$ws.Range("B19").Value = "Thu Aug 21 02:26:43 EDT 2025"
$ws.Range("B20").Value = "Thu Aug 21 02:26:46 EDT 2025"
$ws.Range("A19").Value = "Fail"

$ws = $wb.Worksheets.Item("Void-MRF")
$ws.Range("B2").Value = "Thu Aug 21 02:26:49 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 21 02:26:52 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 21 02:26:54 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 21 02:26:57 EDT 2025"
$ws.Range("B6").Value = "Thu Aug 21 02:26:59 EDT 2025"
$ws.Range("B7").Value = "Thu Aug 21 02:27:02 EDT 2025"
$ws.Range("B8").Value = "Thu Aug 21 02:27:04 EDT 2025"
$ws.Range("B9").Value = "Thu Aug 21 02:27:07 EDT 2025"
$ws.Range("B10").Value = "Thu Aug 21 02:27:09 EDT 2025"
$ws.Range("B11").Value = "Thu Aug 21 02:27:11 EDT 2025"
$ws.Range("A4").Value = "Fail"
$ws.Range("A5").Value = "Fail"
$ws.Range("A7").Value = "Fail"
$ws.Range("A9").Value = "Fail"
$ws.Range("A11").Value = "Fail"

$ws = $wb.Worksheets.Item("Credit-MRF")
$ws.Range("B2").Value = "Thu Aug 21 02:27:14 EDT 2025"
$ws.Range("B3").Value = "Thu Aug 21 02:27:16 EDT 2025"
$ws.Range("B4").Value = "Thu Aug 21 02:27:19 EDT 2025"
$ws.Range("B5").Value = "Thu Aug 21 02:27:21 EDT 2025"
$ws.Range("B6").Value = "Thu Aug 21 02:27:24 EDT 2025"
$ws.Range("B7").Value = "Thu Aug 21 02:27:26 EDT 2025"
$ws.Range("B8").Value = "Thu Aug 21 02:27:28 EDT 2025"
$ws.Range("B9").Value = "Thu Aug 21 02:27:31 EDT 2025"
$ws.Range("B10").Value = "Thu Aug 21 02:27:34 EDT 2025"
$ws.Range("B11").Value = "Thu Aug 21 02:27:36 EDT 2025"
$ws.Range("B12").Value = "Thu Aug 21 02:27:38 EDT 2025"
$ws.Range("B13").Value = "Thu Aug 21 02:27:41 EDT 2025"
$ws.Range("A5").Value = "Fail"
$ws.Range("A6").Value = "Fail"
$ws.Range("A8").Value = "Fail"
$ws.Range("A11").Value = "Fail"

$ws = $wb.Worksheets.Item("Debit-RemID-Pipe")
$ws.Range("B2").Value = "Fri Aug 22 18:07:37 EDT 2025"
$ws.Range("B3").Value = "Fri Aug 22 18:07:41 EDT 2025"
$ws.Range("B4").Value = "Fri Aug 22 18:07:46 EDT 2025"
$ws.Range("B5").Value = "Fri Aug 22 18:07:50 EDT 2025"
$ws.Range("B6").Value = "Fri Aug 22 18:07:55 EDT 2025"
$ws.Range("B7").Value = "Fri Aug 22 18:08:00 EDT 2025"
$ws.Range("B8").Value = "Fri Aug 22 18:08:04 EDT 2025"

$ws = $wb.Worksheets.Item("DebitVoid-RemID-Pipe")
$ws.Range("B2").Value = "Fri Aug 22 18:08:14 EDT 2025"
$ws.Range("B3").Value = "Fri Aug 22 18:08:27 EDT 2025"
$ws.Range("B4").Value = "Fri Aug 22 18:08:36 EDT 2025"
$ws.Range("B5").Value = "Fri Aug 22 18:08:44 EDT 2025"
$ws.Range("B6").Value = "Fri Aug 22 18:08:52 EDT 2025"
$ws.Range("B7").Value = "Fri Aug 22 18:09:00 EDT 2025"
$ws.Range("B8").Value = "Fri Aug 22 18:09:09 EDT 2025"

$ws = $wb.Worksheets.Item("DebitCredit-RemID-Pipe")
$ws.Range("B2").Value = "Fri Aug 22 18:06:34 EDT 2025"
$ws.Range("B3").Value = "Fri Aug 22 18:06:43 EDT 2025"
$ws.Range("B4").Value = "Fri Aug 22 18:06:52 EDT 2025"
$ws.Range("B5").Value = "Fri Aug 22 18:07:01 EDT 2025"
$ws.Range("B6").Value = "Fri Aug 22 18:07:11 EDT 2025"
$ws.Range("B7").Value = "Fri Aug 22 18:07:19 EDT 2025"
$ws.Range("B8").Value = "Fri Aug 22 18:07:28 EDT 2025"
